# Add the new row of data (row 12) to Sheet1, matching the date formatting
# already used in column A, then move the active selection to B15 (no data),
# as captured by the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: 11/10/2025 -> 77 errors
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A12").Value = 45971
$ws.Range("B12").Value = 77

# Move the selection cursor (no associated data) to match the saved view state
$ws.Range("B15").Select()
